$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B30").Value = '\${testBean1}${testBean1}\${testBean2}${testBean2}'
$ws.Range("A30").Value = "EscExprs:"
